$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("K-Fall")
$ws.Range("B2").Value = 81
$ws.Range("B3").Value = 83
$ws.Range("B4").Value = 85
$ws.Range("B5").Value = 87
$ws.Range("B6").Value = 89
$ws.Range("B7").Value = 91
$ws.Range("B8").Value = 93
$ws.Range("B9").Value = 96
$ws.Range("B10").Value = 98
$ws.Range("B11").Value = 100
$ws.Range("B12").Value = 102
$ws.Range("B13").Value = 104
$ws.Range("B14").Value = 106
$ws.Range("B15").Value = 109
$ws.Range("B16").Value = 111
$ws.Range("B21").Value = 122
$ws.Range("B22").Value = 124
$ws.Range("B23").Value = 126

$ws = $wb.Worksheets.Item("K-Spring")
$ws.Range("B2").Value = 74
$ws.Range("B3").Value = 76
$ws.Range("B4").Value = 78
$ws.Range("B5").Value = 81
$ws.Range("B6").Value = 83
$ws.Range("B7").Value = 85
$ws.Range("B8").Value = 87
$ws.Range("B9").Value = 89
$ws.Range("B10").Value = 91
$ws.Range("B11").Value = 94
$ws.Range("B16").Value = 104
$ws.Range("B17").Value = 107
$ws.Range("B18").Value = 109
$ws.Range("B19").Value = 111
$ws.Range("B20").Value = 113
$ws.Range("B21").Value = 115
$ws.Range("B22").Value = 117
$ws.Range("B23").Value = 119
$ws.Range("B24").Value = 122
$ws.Range("B25").Value = 124
$ws.Range("B26").Value = 126
$ws.Range("B27").Value = 128
$ws.Range("B28").Value = 130
$ws.Range("B29").Value = 130

$ws = $wb.Worksheets.Item("1-Fall")
$ws.Range("B2").Value = 68
$ws.Range("B3").Value = 70
$ws.Range("B4").Value = 72
$ws.Range("B5").Value = 74
$ws.Range("B6").Value = 76
$ws.Range("B7").Value = 79
$ws.Range("B8").Value = 81
$ws.Range("B9").Value = 83
$ws.Range("B13").Value = 92
$ws.Range("B14").Value = 94
$ws.Range("B15").Value = 96
$ws.Range("B16").Value = 98
$ws.Range("B17").Value = 100
$ws.Range("B18").Value = 102
$ws.Range("B19").Value = 104
$ws.Range("B20").Value = 107
$ws.Range("B21").Value = 109
$ws.Range("B22").Value = 111
$ws.Range("B23").Value = 113
$ws.Range("B24").Value = 115
$ws.Range("B25").Value = 117
$ws.Range("B26").Value = 120
$ws.Range("B27").Value = 122
$ws.Range("B28").Value = 124
$ws.Range("B29").Value = 126
$ws.Range("B30").Value = 128
$ws.Range("B31").Value = 130

$ws = $wb.Worksheets.Item("1-Spring")
$ws.Range("B2").Value = 81
$ws.Range("B3").Value = 83
$ws.Range("B4").Value = 85
$ws.Range("B5").Value = 87
$ws.Range("B6").Value = 89
$ws.Range("B7").Value = 91
$ws.Range("B8").Value = 93
$ws.Range("B9").Value = 96
$ws.Range("B10").Value = 98
$ws.Range("B11").Value = 100
$ws.Range("B12").Value = 102
$ws.Range("B13").Value = 104
$ws.Range("B14").Value = 106
$ws.Range("B15").Value = 109
$ws.Range("B16").Value = 111
$ws.Range("B17").Value = 113
$ws.Range("B18").Value = 115
$ws.Range("B19").Value = 117
$ws.Range("B20").Value = 119
$ws.Range("B21").Value = 122
$ws.Range("B22").Value = 124
$ws.Range("B23").Value = 126
$ws.Range("B24").Value = 128
$ws.Range("B25").Value = 130
$ws.Range("B26").Value = 130
$ws.Range("B27").Value = 130
$ws.Range("B28").Value = 130
$ws.Range("B29").Value = 130
$ws.Range("B30").Value = 130
$ws.Range("B31").Value = 130

$ws = $wb.Worksheets.Item("2-Fall")
$ws.Range("B2").Value = 55
$ws.Range("B3").Value = 57
$ws.Range("B4").Value = 59
$ws.Range("B5").Value = 62
$ws.Range("B6").Value = 64
$ws.Range("B7").Value = 66
$ws.Range("B8").Value = 68
$ws.Range("B11").Value = 75
$ws.Range("B12").Value = 77
$ws.Range("B13").Value = 79
$ws.Range("B14").Value = 81
$ws.Range("B15").Value = 83
$ws.Range("B16").Value = 85
$ws.Range("B17").Value = 87
$ws.Range("B18").Value = 90
$ws.Range("B19").Value = 92
$ws.Range("B20").Value = 94
$ws.Range("B21").Value = 96
$ws.Range("B22").Value = 98
$ws.Range("B23").Value = 100
$ws.Range("B24").Value = 103
$ws.Range("B25").Value = 105
$ws.Range("B26").Value = 107
$ws.Range("B27").Value = 109
$ws.Range("B28").Value = 111
$ws.Range("B31").Value = 118

$ws = $wb.Worksheets.Item("2-Spring")
$ws.Range("B2").Value = 49
$ws.Range("B3").Value = 51
$ws.Range("B4").Value = 53
$ws.Range("B5").Value = 55
$ws.Range("B6").Value = 57
$ws.Range("B7").Value = 60
$ws.Range("B8").Value = 62
$ws.Range("B9").Value = 64
$ws.Range("B10").Value = 66
$ws.Range("B11").Value = 68
$ws.Range("B14").Value = 75
$ws.Range("B15").Value = 77
$ws.Range("B16").Value = 79
$ws.Range("B17").Value = 81
$ws.Range("B18").Value = 83
$ws.Range("B19").Value = 85
$ws.Range("B20").Value = 88
$ws.Range("B21").Value = 90
$ws.Range("B22").Value = 92
$ws.Range("B23").Value = 94
$ws.Range("B24").Value = 96
$ws.Range("B28").Value = 105
$ws.Range("B29").Value = 107
$ws.Range("B30").Value = 109
$ws.Range("B31").Value = 111
